# Edit the abstract text per the commit diff:
#  1. "reinstated" -> "resumed"
#  2. Restructure the closing two sentences about monitoring data / dam
#     operations and the importance of understanding river conditions.

$d = $word.ActiveDocument

# 1) "DWR reinstated the lower Feather River" -> "DWR resumed the lower Feather River"
$d.Content.Find.Execute(
    "DWR reinstated the lower Feather River",
    $true, $true, $false, $false, $false,
    $true, 1, $false,
    "DWR resumed the lower Feather River",
    2
)

# 2) Rework the final two sentences of the paragraph.
$old2 = "Data from this monitoring will also inform effective management of dam operations and implementation of any associated restoration activities. Having a thorough understanding of how river conditions and habitats affect the distribution, abundance, and behavior of downstream fish communities is crucial."
$new2 = "Due to the importance of having a thorough understanding of how river conditions and habitats affect the distribution, abundance, and behavior of downstream fish communities, data from this monitoring will inform effective management of dam operations and implementation of any associated restoration activities."

$d.Content.Find.Execute(
    $old2,
    $true, $true, $false, $false, $false,
    $true, 1, $false,
    $new2,
    2
)
